$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update labels and values (rows 1-3)
$ws.Range("A1").Value = "Dia da cotação: "
$ws.Range("B1").Value = "Tue May 21 2024 14:51:26 GMT-0400 (Amazon Standard Time)"

$ws.Range("A2").Value = "Cotação dolar: "
$ws.Range("B2").Value = "1"

$ws.Range("A3").Value = "Cotação real: "
$ws.Range("B3").Value = "5.12"

# Add new row 4
$ws.Range("A4").Value = "100 dolar convertido: "
$ws.Range("B4").Value = "512 em real"

# Update column widths for columns A-D, and remove custom width from column E
$ws.Range("A1:D1").ColumnWidth = 30.83203125
$ws.Columns.Item(5).ColumnWidth = 8.43
